$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1-6. Remove spell-check proofErr run-splits by doing a Find & Replace of the
#      full sentence text; Word's COM layer re-merges the runs into one and
#      drops the proofErr markers automatically when the replacement text is
#      written back.
# ---------------------------------------------------------------------------

[void]$d.Content.Find.Execute(
    "husk å invitere Eirik, Asle og Atle i Skrum ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "husk å invitere Eirik, Asle og Atle i Skrum ", 2)

[void]$d.Content.Find.Execute(
    "Tilsendt dokument (docs) – se igjennom det for å få innblikk ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Tilsendt dokument (docs) – se igjennom det for å få innblikk ", 2)

[void]$d.Content.Find.Execute(
    "Få litt mer beskrivelse fra Eirik om bedriftens forventninger og forvetningsverdi. ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Få litt mer beskrivelse fra Eirik om bedriftens forventninger og forvetningsverdi. ", 2)

[void]$d.Content.Find.Execute(
    "Finne flere forskningpaper som vi kan benytte", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Finne flere forskningpaper som vi kan benytte", 2)

[void]$d.Content.Find.Execute(
    "- Teste IE8+ (senere i sprint) ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "- Teste IE8+ (senere i sprint) ", 2)

[void]$d.Content.Find.Execute(
    " Ta i bruk ifb forskning og dokumentasjon", $true, $false, $false,
    $false, $false, $true, 1, $false,
    " Ta i bruk ifb forskning og dokumentasjon", 2)

# ---------------------------------------------------------------------------
# 7. Remove the old "_GoBack" bookmark (it sat at the end of the "Finne flere
#    forskningpaper..." paragraph) - Word re-creates it at the location of
#    the newest edit, which is the new paragraph inserted below.
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 8. Insert a brand-new paragraph right after the
#    "- informasjonsevaluering ... Ta i bruk ifb forskning og dokumentasjon"
#    paragraph, containing the new sentence about the RIS technique, and
#    re-create the "_GoBack" bookmark collapsed right after the new text.
# ---------------------------------------------------------------------------

$anchorRange = $d.Content
[void]$anchorRange.Find.Execute(
    " Ta i bruk ifb forskning og dokumentasjon", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0)

$anchorIndex = $anchorRange.Paragraphs(1).Index
[void]$anchorRange.InsertParagraphAfter()

$newIndex = $anchorIndex + 1
$newParaRange = $d.Paragraphs($newIndex).Range

# Write the sentence with a one-character sentinel on the end; this lets us
# collapse a Find-match range to a position that is NOT the very last
# character slot of the paragraph (collapsing right at that slot mis-anchors
# newly-added bookmarks), add the bookmark, then strip the sentinel.
$sentinelText = "- Bruke RIS-teknikken/modell -> for å lage grafiske modeller av en prosess (MUST)Z"
$newParaRange.Text = $sentinelText

$sentinelRange = $d.Paragraphs($newIndex).Range
[void]$sentinelRange.Find.Execute("Z", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sentinelRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $sentinelRange)

$deleteRange = $d.Paragraphs($newIndex).Range
[void]$deleteRange.Find.Execute("Z", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deleteRange.Text = ""
